$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A12").Value = "Tempo"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "analogue"
$ws.Range("D12").Value = "C4"
Write-Output ($ws.Cells.Item(12,1).Value2)
Write-Output ($ws.Cells.Item(12,2).Value2)
Write-Output ($ws.Cells.Item(12,3).Value2)
Write-Output ($ws.Cells.Item(12,4).Value2)
